$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 19374.875
$ws.Range("I70").Value = 4000
$ws.Range("J70").Value = 21571.285
$ws.Range("K70").Value = 12000
$ws.Range("L70").Value = 64713.855
$ws.Range("M70").Value = -11730
$ws.Range("N70").Value = -65253.855
$ws.Range("H73").Value = 19374.875
$ws.Range("I73").Value = 4000
$ws.Range("J73").Value = 21571.285
$ws.Range("K73").Value = 12000
$ws.Range("L73").Value = 64713.855
$ws.Range("M73").Value = -11064
$ws.Range("N73").Value = -66585.855
$ws.Range("H80").Value = 16798.334
$ws.Range("I80").Value = 197.5
$ws.Range("J80").Value = 50000
$ws.Range("K80").Value = 592.5
$ws.Range("L80").Value = 150000
$ws.Range("M80").Value = 405.5
$ws.Range("H83").Value = 16798.334
$ws.Range("I83").Value = 197.5
$ws.Range("J83").Value = 50000
$ws.Range("K83").Value = 1777.5
$ws.Range("L83").Value = 450000
$ws.Range("M83").Value = 3214.5
$ws.Range("H140").Value = 99995
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 99995
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 99995
$ws.Range("N140").Value = -110355

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 2005
$ws.Range("I3").Value = 2005
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 2005
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -1890
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = $null
$ws.Range("N5").Value = $null
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = $null
$ws.Range("H45").Value = 1437.5
$ws.Range("I45").Value = 1437.5
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 1437.5
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -1060.5
$ws.Range("H61").Value = 2849.8333
$ws.Range("I61").Value = 2849.8333
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2849.8333
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2637.8333
$ws.Range("H88").Value = 2048.9
$ws.Range("I88").Value = 1428.4286
$ws.Range("J88").Value = 3496.6667
$ws.Range("K88").Value = 1428.4286
$ws.Range("L88").Value = 3496.6667
$ws.Range("M88").Value = -1022.4286
$ws.Range("N88").Value = -4308.6667
$ws.Range("H91").Value = 2048.9
$ws.Range("I91").Value = 1428.4286
$ws.Range("J91").Value = 3496.6667
$ws.Range("K91").Value = 1428.4286
$ws.Range("L91").Value = 3496.6667
$ws.Range("M91").Value = -24.42859999999996
$ws.Range("N91").Value = -6304.6667
$ws.Range("H136").Value = 2849.8333
$ws.Range("I136").Value = 2849.8333
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 8549.499899999999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -5999.499899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = $null
$ws.Range("N4").Value = $null
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = $null
$ws.Range("H11").Value = 298.5
$ws.Range("I11").Value = 331.33334
$ws.Range("J11").Value = 200
$ws.Range("K11").Value = 331.33334
$ws.Range("L11").Value = 200
$ws.Range("M11").Value = -191.33334
$ws.Range("N11").Value = -480
$ws.Range("H12").Value = 9999
$ws.Range("I12").Value = 9999
$ws.Range("J12").Value = 9999
$ws.Range("K12").Value = 9999
$ws.Range("L12").Value = 9999
$ws.Range("M12").Value = -9831
$ws.Range("N12").Value = -10335
$ws.Range("H22").Value = 849.44446
$ws.Range("I22").Value = 849.44446
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 849.44446
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -676.44446
$ws.Range("H25").Value = 3000
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 3000
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 3000
$ws.Range("N25").Value = -3470
$ws.Range("H29").Value = 2555.5
$ws.Range("I29").Value = 111
$ws.Range("J29").Value = 5000
$ws.Range("K29").Value = 111
$ws.Range("L29").Value = 5000
$ws.Range("M29").Value = 178
$ws.Range("N29").Value = -5578
$ws.Range("H86").Value = 1499.5
$ws.Range("I86").Value = 1499.5
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 1499.5
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -376.5
$ws.Range("H89").Value = 1499.5
$ws.Range("I89").Value = 1499.5
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 7497.5
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -1881.5
$ws.Range("H99").Value = 1000
$ws.Range("I99").Value = 1000
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1000
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 498
$ws.Range("N99").Value = $null
$ws.Range("H108").Value = 91498.336
$ws.Range("I108").Value = 75000
$ws.Range("J108").Value = 99747.5
$ws.Range("K108").Value = 75000
$ws.Range("L108").Value = 99747.5
$ws.Range("M108").Value = -71160
$ws.Range("N108").Value = -107427.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 3459.6
$ws.Range("I15").Value = 300
$ws.Range("J15").Value = 4249.5
$ws.Range("K15").Value = 300
$ws.Range("L15").Value = 4249.5
$ws.Range("M15").Value = -130
$ws.Range("N15").Value = -4589.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H43").Value = 2000
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 2000
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 6000
$ws.Range("N43").Value = -6228
$ws.Range("H80").Value = 5866.6665
$ws.Range("I80").Value = 4400
$ws.Range("J80").Value = 6285.7144
$ws.Range("K80").Value = 13200
$ws.Range("L80").Value = 18857.1432
$ws.Range("M80").Value = -12264
$ws.Range("N80").Value = -20729.1432
$ws.Range("H83").Value = 5866.6665
$ws.Range("I83").Value = 4400
$ws.Range("J83").Value = 6285.7144
$ws.Range("K83").Value = 39600
$ws.Range("L83").Value = 56571.4296
$ws.Range("M83").Value = -34920
$ws.Range("N83").Value = -65931.4296
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = $null

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 131.5
$ws.Range("I2").Value = 123.84615
$ws.Range("J2").Value = 164.66667
$ws.Range("K2").Value = 123.84615
$ws.Range("L2").Value = 164.66667
$ws.Range("M2").Value = -10.84614999999999
$ws.Range("H47").Value = 10000
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 10000
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 10000
$ws.Range("N47").Value = -11136
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = $null
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = $null
$ws.Range("H107").Value = 434.625
$ws.Range("I107").Value = 375.25
$ws.Range("J107").Value = 494
$ws.Range("K107").Value = 375.25
$ws.Range("L107").Value = 494
$ws.Range("M107").Value = 1544.75
$ws.Range("N107").Value = -4334
$ws.Range("H113").Value = 1500
$ws.Range("I113").Value = 1500
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1500
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 670
$ws.Range("H132").Value = 2280.6667
$ws.Range("I132").Value = 2236.8
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 6710.400000000001
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -4180.400000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5927.317
$ws.Range("I46").Value = 4000
$ws.Range("J46").Value = 6026.154
$ws.Range("K46").Value = 4000
$ws.Range("L46").Value = 6026.154
$ws.Range("M46").Value = -3812
$ws.Range("H74").Value = 10000
$ws.Range("I74").Value = 10000
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 10000
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -9002
$ws.Range("H77").Value = 10000
$ws.Range("I77").Value = 10000
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 30000
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -25008
$ws.Range("H100").Value = 3730.1428
$ws.Range("I100").Value = 4185.3335
$ws.Range("J100").Value = 999
$ws.Range("K100").Value = 4185.3335
$ws.Range("L100").Value = 999
$ws.Range("M100").Value = -3644.3335
$ws.Range("N100").Value = -2081

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 1015243.8
$ws.Range("I24").Value = 2505504.5
$ws.Range("J24").Value = 21736.666
$ws.Range("K24").Value = 2505504.5
$ws.Range("L24").Value = 21736.666
$ws.Range("M24").Value = -2505274.5
$ws.Range("H28").Value = 15258.5
$ws.Range("I28").Value = 15258.5
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 15258.5
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -14910.5
$ws.Range("H31").Value = 20017
$ws.Range("I31").Value = 20017
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 20017
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -19669
$ws.Range("H110").Value = 44500
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 44500
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 44500
$ws.Range("N110").Value = -52680
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").Value = $null
